$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - previously only "no"; now coded as "yes" with full classification
$ws.Range("B7").Value = "yes"
$ws.Range("C7").Value = "action"
$ws.Range("D7").Value = "measures"
$ws.Range("E7").Value = "global"
$ws.Range("F7").Value = "distant future"
$ws.Range("G7").Value = "egalitarian"
$ws.Range("H7").Value = "Judgement on the need of a committed and global effort.  "

# Row 12 - "other(needs)" / "other(resources)" simplified to "needs" / "resources"
$ws.Range("C12").Value = "needs"
$ws.Range("D12").Value = "resources"

# Row 13 - Time updated from "nearby future" to "present"
$ws.Range("F13").Value = "present"

# Row 15 - no longer relevant; clear classification, relevance flips to "no"
$ws.Range("B15").Value = "no"
$ws.Range("C15:H15").ClearContents()

# Row 17 - resources note expanded to include technological resources
$ws.Range("D17").Value = "financial resources, technological resources"

# Row 27 - no longer relevant; clear classification, relevance flips to "no"
$ws.Range("B27").Value = "no"
$ws.Range("C27:H27").ClearContents()

# Row 30 - previously only "no"; now coded as "yes" with full classification
$ws.Range("B30").Value = "yes"
$ws.Range("C30").Value = "loss and damage, new UNFCCC policy"
$ws.Range("D30").Value = "n.a."
$ws.Range("E30").Value = "n.a."
$ws.Range("F30").Value = "present"
$ws.Range("G30").Value = "general normative statement"
$ws.Range("H30").Value = "General call for the inclusion of loss and damage in the new agreement. "

# Row 36 - previously only "no"; now coded as "yes" with full classification
$ws.Range("B36").Value = "yes"
$ws.Range("C36").Value = "action"
$ws.Range("D36").Value = "n.a."
$ws.Range("E36").Value = "n.a."
$ws.Range("F36").Value = "present"
$ws.Range("G36").Value = "utilitarian"
$ws.Range("H36").Value = "Presenting the value of needing to protect the survival of humankind, thus benefit of all. "

# Selection / scroll position update
$ws.Range("C7").Select()
